$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 50500
$ws.Range("J3").Value = 50500
$ws.Range("L3").Value = 50500
$ws.Range("N3").Value = -50728

$ws.Range("H12").Value = 537.1667
$ws.Range("I12").Value = 558.6
$ws.Range("J12").Value = 430
$ws.Range("K12").Value = 558.6
$ws.Range("L12").Value = 430
$ws.Range("M12").Value = -388.6
$ws.Range("N12").Value = -770

$ws.Range("H19").Value = 1685
$ws.Range("I19").Value = 1885
$ws.Range("K19").Value = 1885
$ws.Range("M19").Value = -1710

$ws.Range("H51").Value = 17084.9
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 17872.223
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 17872.223
$ws.Range("M51").Value = -9515
$ws.Range("N51").Value = -18840.223

$ws.Range("H62").Value = 97439176
$ws.Range("I62").Value = 140743490
$ws.Range("K62").Value = 140743490
$ws.Range("M62").Value = -140742866

$ws.Range("H65").Value = 97439176
$ws.Range("I65").Value = 140743490
$ws.Range("K65").Value = 703717450
$ws.Range("M65").Value = -703714330

$ws.Range("H100").Value = 22596.213
$ws.Range("J100").Value = 9440
$ws.Range("L100").Value = 9440
$ws.Range("N100").Value = -10522

$ws.Range("H102").Value = 50500
$ws.Range("J102").Value = 50500
$ws.Range("L102").Value = 50500
$ws.Range("N102").Value = -56990

$ws.Range("H107").Value = 16129852
$ws.Range("I107").Value = 17857930
$ws.Range("K107").Value = 17857930
$ws.Range("M107").Value = -17856010

$ws.Range("H132").Value = 1717.25
$ws.Range("I132").Value = 1534
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4602
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2072
$ws.Range("N132").Value = -14060

$ws.Range("H138").Value = 2865.9565
$ws.Range("I138").Value = 2065.125
$ws.Range("J138").Value = 3034.5527
$ws.Range("K138").Value = 6195.375
$ws.Range("L138").Value = 9103.658100000001
$ws.Range("M138").Value = -1055.375
$ws.Range("N138").Value = -19383.6581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0

$ws.Range("H61").Value = 6863.4546
$ws.Range("I61").Value = 4599.6
$ws.Range("J61").Value = 8750
$ws.Range("K61").Value = 4599.6
$ws.Range("L61").Value = 8750
$ws.Range("M61").Value = -4387.6
$ws.Range("N61").Value = -9174

$ws.Range("H74").Value = 591916.3
$ws.Range("I74").Value = 1429368.1
$ws.Range("J74").Value = 5700
$ws.Range("K74").Value = 1429368.1
$ws.Range("L74").Value = 5700
$ws.Range("M74").Value = -1428494.1
$ws.Range("N74").Value = -7448

$ws.Range("H77").Value = 591916.3
$ws.Range("I77").Value = 1429368.1
$ws.Range("J77").Value = 5700
$ws.Range("K77").Value = 7146840.5
$ws.Range("L77").Value = 28500
$ws.Range("M77").Value = -7142472.5
$ws.Range("N77").Value = -37236

$ws.Range("H106").Value = 81342.5
$ws.Range("J106").Value = 81342.5
$ws.Range("L106").Value = 81342.5
$ws.Range("N106").Value = -83866.5

$ws.Range("H132").Value = 6183.6665
$ws.Range("I132").Value = 4515.5454
$ws.Range("K132").Value = 13546.6362
$ws.Range("M132").Value = -11016.6362

$ws.Range("H136").Value = 6863.4546
$ws.Range("I136").Value = 4599.6
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 13798.8
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -11248.8
$ws.Range("N136").Value = -31350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 103706.2
$ws.Range("I86").Value = 4009.0588
$ws.Range("J86").Value = 668656.7
$ws.Range("K86").Value = 4009.0588
$ws.Range("L86").Value = 668656.7
$ws.Range("M86").Value = -2886.0588
$ws.Range("N86").Value = -670902.7

$ws.Range("H89").Value = 103706.2
$ws.Range("I89").Value = 4009.0588
$ws.Range("J89").Value = 668656.7
$ws.Range("K89").Value = 20045.294
$ws.Range("L89").Value = 3343283.5
$ws.Range("M89").Value = -14429.294
$ws.Range("N89").Value = -3354515.5

$ws.Range("H134").Value = 2724.2258
$ws.Range("I134").Value = 1117.52
$ws.Range("K134").Value = 3352.56
$ws.Range("M134").Value = -817.5599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4076.2
$ws.Range("I16").Value = 3860.2144
$ws.Range("K16").Value = 3860.2144
$ws.Range("M16").Value = -3573.2144

$ws.Range("H31").Value = 111118376
$ws.Range("I31").Value = 500001180
$ws.Range("J31").Value = 9000
$ws.Range("K31").Value = 500001180
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = -500000885
$ws.Range("N31").Value = -9590

$ws.Range("H34").Value = 111118376
$ws.Range("I34").Value = 500001180
$ws.Range("J34").Value = 9000
$ws.Range("K34").Value = 500001180
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -500000978
$ws.Range("N34").Value = -9404

$ws.Range("H43").Value = 63458.305
$ws.Range("J43").Value = 63458.305
$ws.Range("L43").Value = 63458.305
$ws.Range("N43").Value = -63826.305

$ws.Range("H96").Value = 33434.4
$ws.Range("J96").Value = 33434.4
$ws.Range("L96").Value = 33434.4
$ws.Range("N96").Value = -38926.4

$ws.Range("H101").Value = 63458.305
$ws.Range("J101").Value = 63458.305
$ws.Range("L101").Value = 63458.305
$ws.Range("N101").Value = -69948.30499999999

$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253

$ws.Range("H113").Value = 4076.2
$ws.Range("I113").Value = 3860.2144
$ws.Range("K113").Value = 3860.2144
$ws.Range("M113").Value = -1690.2144

$ws.Range("H134").Value = 7989
$ws.Range("I134").Value = 8004.72
$ws.Range("J134").Value = 7923.5
$ws.Range("K134").Value = 24014.16
$ws.Range("L134").Value = 23770.5
$ws.Range("M134").Value = -21479.16
$ws.Range("N134").Value = -28840.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2999.5
$ws.Range("J2").Value = 4442.25
$ws.Range("L2").Value = 26653.5
$ws.Range("N2").Value = -26879.5

$ws.Range("H38").Value = 118.454544
$ws.Range("I38").Value = 19.25
$ws.Range("K38").Value = 57.75
$ws.Range("M38").Value = 289.25

$ws.Range("H47").Value = 100502.4
$ws.Range("I47").Value = 111447.11
$ws.Range("K47").Value = 334341.33
$ws.Range("M47").Value = -333910.33

$ws.Range("H74").Value = 9128.25
$ws.Range("J74").Value = 9500
$ws.Range("L74").Value = 28500
$ws.Range("N74").Value = -30622

$ws.Range("H77").Value = 9128.25
$ws.Range("J77").Value = 9500
$ws.Range("L77").Value = 85500
$ws.Range("N77").Value = -96108

$ws.Range("H80").Value = 4399.6665
$ws.Range("J80").Value = 4399.6665
$ws.Range("L80").Value = 13198.9995
$ws.Range("N80").Value = -15070.9995

$ws.Range("H83").Value = 4399.6665
$ws.Range("J83").Value = 4399.6665
$ws.Range("L83").Value = 39596.9985
$ws.Range("N83").Value = -48956.9985

$ws.Range("H92").Value = 1116.7
$ws.Range("I92").Value = 1220.9231
$ws.Range("K92").Value = 3662.7693
$ws.Range("M92").Value = -2414.7693

$ws.Range("H114").Value = 3266.4
$ws.Range("I114").Value = 325.25
$ws.Range("K114").Value = 975.75
$ws.Range("M114").Value = 2278.25

$ws.Range("H131").Value = 16671543
$ws.Range("I131").Value = 37037940
$ws.Range("J131").Value = 8128.273
$ws.Range("K131").Value = 111113820
$ws.Range("L131").Value = 24384.819
$ws.Range("M131").Value = -111108780
$ws.Range("N131").Value = -34464.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3132.6
$ws.Range("J126").Value = 4814
$ws.Range("L126").Value = 14442
$ws.Range("N126").Value = -19382

$ws.Range("H129").Value = 92780
$ws.Range("J129").Value = 92780
$ws.Range("L129").Value = 92780

$ws.Range("H132").Value = 3857.2727
$ws.Range("I132").Value = 2380.353
$ws.Range("J132").Value = 5426.5
$ws.Range("K132").Value = 7141.059
$ws.Range("L132").Value = 16279.5
$ws.Range("M132").Value = -4611.059
$ws.Range("N132").Value = -21339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2282.7144
$ws.Range("I82").Value = 1950.8889
$ws.Range("J82").Value = 2880
$ws.Range("K82").Value = 1950.8889
$ws.Range("L82").Value = 2880
$ws.Range("M82").Value = -1589.8889
$ws.Range("N82").Value = -3602

$ws.Range("H85").Value = 2282.7144
$ws.Range("I85").Value = 1950.8889
$ws.Range("J85").Value = 2880
$ws.Range("K85").Value = 1950.8889
$ws.Range("L85").Value = 2880
$ws.Range("M85").Value = -702.8888999999999
$ws.Range("N85").Value = -5376

$ws.Range("H93").Value = 3808.88
$ws.Range("I93").Value = 2800.7334
$ws.Range("K93").Value = 2800.7334
$ws.Range("M93").Value = -1552.7334

$ws.Range("H130").Value = 26057.25
$ws.Range("J130").Value = 26057.25
$ws.Range("L130").Value = 26057.25
$ws.Range("N130").Value = -36097.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4149.3076
$ws.Range("I122").Value = 4120.0835
$ws.Range("K122").Value = 12360.2505
$ws.Range("M122").Value = -9910.250499999998

$ws.Range("H132").Value = 5389.5264
$ws.Range("I132").Value = 3160.75
$ws.Range("J132").Value = 7865.9443
$ws.Range("K132").Value = 9482.25
$ws.Range("L132").Value = 23597.8329
$ws.Range("M132").Value = -6952.25
$ws.Range("N132").Value = -28657.8329

# Special case: ARM row 57 - M57 cell removed entirely (was -8016)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M57").ClearContents()

# Special case: GSM row 129 - N129 cell newly added (-102780)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N129").Value = -102780
